# Auto-generated Excel COM-interop script applying the numeric updates
# described in the commit diff for Sargatanas_Profits.xlsx (per-sheet Leve profit data).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 59948556
$ws.Range("I86").Value = 87953450
$ws.Range("J86").Value = 7939475.5
$ws.Range("K86").Value = 87953450
$ws.Range("L86").Value = 7939475.5
$ws.Range("M86").Value = -87952327
$ws.Range("N86").Value = -7941721.5

$ws.Range("H89").Value = 59948556
$ws.Range("I89").Value = 87953450
$ws.Range("J89").Value = 7939475.5
$ws.Range("K89").Value = 439767250
$ws.Range("L89").Value = 39697377.5
$ws.Range("M89").Value = -439761634
$ws.Range("N89").Value = -39708609.5

$ws.Range("H116").Value = 41670800
$ws.Range("I116").Value = 83335830
$ws.Range("K116").Value = 83335830
$ws.Range("M116").Value = -83332388

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 825.3469
$ws.Range("I132").Value = 647.9149
$ws.Range("K132").Value = 1943.7447
$ws.Range("M132").Value = 586.2553

$ws.Range("H137").Value = 2358.5454
$ws.Range("I137").Value = 1992.4286
$ws.Range("K137").Value = 5977.2858
$ws.Range("M137").Value = -3427.2858

$ws.Range("H138").Value = 6348.65
$ws.Range("J138").Value = 6472.971
$ws.Range("L138").Value = 19418.913
$ws.Range("N138").Value = -29698.913

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4019070.5
$ws.Range("I32").Value = 4100378.2
$ws.Range("K32").Value = 4100378.2
$ws.Range("M32").Value = -4100091.2

$ws.Range("H45").Value = 3666.4443
$ws.Range("I45").Value = 2999.6667
$ws.Range("J45").Value = 3999.8333
$ws.Range("K45").Value = 2999.6667
$ws.Range("L45").Value = 3999.8333
$ws.Range("M45").Value = -2622.6667
$ws.Range("N45").Value = -4753.8333

$ws.Range("H61").Value = 28579326
$ws.Range("I61").Value = 4530
$ws.Range("K61").Value = 4530
$ws.Range("M61").Value = -4318

$ws.Range("H74").Value = 41839.08
$ws.Range("I74").Value = 64801.062
$ws.Range("K74").Value = 64801.062
$ws.Range("M74").Value = -63927.062

$ws.Range("H77").Value = 41839.08
$ws.Range("I77").Value = 64801.062
$ws.Range("K77").Value = 324005.31
$ws.Range("M77").Value = -319637.31

$ws.Range("H122").Value = 5639.7
$ws.Range("I122").Value = 5377.4443
$ws.Range("K122").Value = 16132.3329
$ws.Range("M122").Value = -13682.3329

$ws.Range("H123").Value = 80000
$ws.Range("J123").Value = 80000
$ws.Range("L123").Value = 80000
$ws.Range("N123").Value = -89800

$ws.Range("H136").Value = 28579326
$ws.Range("I136").Value = 4530
$ws.Range("K136").Value = 13590
$ws.Range("M136").Value = -11040

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6949529.5
$ws.Range("I20").Value = 9263082
$ws.Range("K20").Value = 9263082
$ws.Range("M20").Value = -9262835

$ws.Range("H107").Value = 30406580
$ws.Range("J107").Value = 1633
$ws.Range("L107").Value = 1633
$ws.Range("N107").Value = -5473

$ws.Range("H134").Value = 5325376
$ws.Range("I134").Value = 10872490
$ws.Range("K134").Value = 32617470
$ws.Range("M134").Value = -32614935

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10914.075
$ws.Range("I31").Value = 5005.5625
$ws.Range("J31").Value = 14853.083
$ws.Range("K31").Value = 5005.5625
$ws.Range("L31").Value = 14853.083
$ws.Range("M31").Value = -4710.5625
$ws.Range("N31").Value = -15443.083

$ws.Range("H34").Value = 10914.075
$ws.Range("I34").Value = 5005.5625
$ws.Range("J34").Value = 14853.083
$ws.Range("K34").Value = 5005.5625
$ws.Range("L34").Value = 14853.083
$ws.Range("M34").Value = -4803.5625
$ws.Range("N34").Value = -15257.083

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H51").Value = 45633.332
$ws.Range("J51").Value = 50000
$ws.Range("L51").Value = 50000
$ws.Range("N51").Value = -51472

$ws.Range("H52").Value = 100260
$ws.Range("I52").Value = 80000
$ws.Range("J52").Value = 110390
$ws.Range("K52").Value = 80000
$ws.Range("L52").Value = 110390
$ws.Range("M52").Value = -79706
$ws.Range("N52").Value = -110978

$ws.Range("H59").Value = 99772.25
$ws.Range("J59").Value = 99772.25
$ws.Range("L59").Value = 99772.25
$ws.Range("N59").Value = -102062.25

$ws.Range("H60").Value = 19666.666
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H61").Value = 45633.332
$ws.Range("J61").Value = 50000
$ws.Range("L61").Value = 50000
$ws.Range("N61").Value = -50696

$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

$ws.Range("H100").Value = 46346
$ws.Range("J100").Value = 46346
$ws.Range("L100").Value = 46346
$ws.Range("N100").Value = -48510

$ws.Range("H105").Value = 4472504
$ws.Range("I105").Value = 8942134
$ws.Range("K105").Value = 8942134
$ws.Range("M105").Value = -8940387

$ws.Range("H106").Value = 50118
$ws.Range("J106").Value = 50118
$ws.Range("L106").Value = 50118
$ws.Range("N106").Value = -52642

$ws.Range("H132").Value = 7843.2583
$ws.Range("I132").Value = 5721.1665
$ws.Range("K132").Value = 17163.4995
$ws.Range("M132").Value = -14633.4995

$ws.Range("H134").Value = 7095.2144
$ws.Range("I134").Value = 2540.9
$ws.Range("K134").Value = 7622.700000000001
$ws.Range("M134").Value = -5087.700000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 166668930
$ws.Range("J117").Value = 200002990
$ws.Range("L117").Value = 600008970
$ws.Range("N117").Value = -600015854

$ws.Range("H120").Value = 20305.23
$ws.Range("J120").Value = 24099.3
$ws.Range("L120").Value = 72297.89999999999
$ws.Range("N120").Value = -81973.89999999999

$ws.Range("H131").Value = 59032.11
$ws.Range("J131").Value = 75060
$ws.Range("L131").Value = 225180
$ws.Range("N131").Value = -235260

$ws.Range("H137").Value = 227825.44
$ws.Range("I137").Value = 503649.5
$ws.Range("J137").Value = 149018.58
$ws.Range("K137").Value = 1510948.5
$ws.Range("L137").Value = 447055.74
$ws.Range("M137").Value = -1505848.5
$ws.Range("N137").Value = -457255.74

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3054.12
$ws.Range("I102").Value = 2820.0227
$ws.Range("J102").Value = 4770.8335
$ws.Range("K102").Value = 2820.0227
$ws.Range("L102").Value = 4770.8335
$ws.Range("M102").Value = -1198.0227
$ws.Range("N102").Value = -8014.8335

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4260.65
$ws.Range("I7").Value = 3479.077
$ws.Range("J7").Value = 5712.143
$ws.Range("K7").Value = 3479.077
$ws.Range("L7").Value = 5712.143
$ws.Range("M7").Value = -3367.077
$ws.Range("N7").Value = -5936.143

$ws.Range("H122").Value = 6654.2964
$ws.Range("I122").Value = 5743.091
$ws.Range("K122").Value = 17229.273
$ws.Range("M122").Value = -14779.273

$ws.Range("H126").Value = 4260.65
$ws.Range("I126").Value = 3479.077
$ws.Range("J126").Value = 5712.143
$ws.Range("K126").Value = 10437.231
$ws.Range("L126").Value = 17136.429
$ws.Range("M126").Value = -7967.231
$ws.Range("N126").Value = -22076.429

$ws.Range("H132").Value = 16137751
$ws.Range("I132").Value = 26320752
$ws.Range("J132").Value = 14666.583
$ws.Range("K132").Value = 78962256
$ws.Range("L132").Value = 43999.749
$ws.Range("M132").Value = -78959726
$ws.Range("N132").Value = -49059.749

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 153998.75
$ws.Range("I62").Value = 300000
$ws.Range("K62").Value = 300000
$ws.Range("M62").Value = -299376

$ws.Range("H65").Value = 153998.75
$ws.Range("I65").Value = 300000
$ws.Range("K65").Value = 1500000
$ws.Range("M65").Value = -1496880

$ws.Range("H81").Value = 17559910
$ws.Range("I81").Value = 1168490.4
$ws.Range("K81").Value = 2336980.8
$ws.Range("M81").Value = -2335919.8

$ws.Range("H84").Value = 17559910
$ws.Range("I84").Value = 1168490.4
$ws.Range("K84").Value = 11684904
$ws.Range("M84").Value = -11679600

$ws.Range("H104").Value = 22750.834
$ws.Range("J104").Value = 22750.834
$ws.Range("L104").Value = 22750.834
$ws.Range("N104").Value = -29738.834

$ws.Range("H122").Value = 204544.95
$ws.Range("I122").Value = 238876.4
$ws.Range("K122").Value = 716629.2
$ws.Range("M122").Value = -714179.2

$ws.Range("H126").Value = 1599.3
$ws.Range("I126").Value = 1770.4286
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 5311.2858
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -2841.2858
$ws.Range("N126").Value = -8540

$ws.Range("H132").Value = 9365.405000000001
$ws.Range("I132").Value = 9977.885
$ws.Range("J132").Value = 7917.727
$ws.Range("K132").Value = 29933.655
$ws.Range("L132").Value = 23753.181
$ws.Range("M132").Value = -27403.655
$ws.Range("N132").Value = -28813.181
